$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")
$ws.Activate()

# userid -> profile image file name, in the same order as the Users sheet
$users = @(
    @("a_wilson", "a_wilson.jpg"),
    @("j_thomas", "j_thomas.jpg"),
    @("p_harris", "p_harris.jpg"),
    @("j_green", "j_green.jpg"),
    @("s_bowers", "s_bowers.jpg"),
    @("l_denton", "l_denton.jpg"),
    @("a_foster", "a_foster.jpg"),
    @("e_hanson", "e_hanson.jpg"),
    @("l_johnson", "l_johnson.jpg"),
    @("e_reese", "e_reese.jpg")
)

$row = 67
foreach ($u in $users) {
    $ws.Cells.Item($row, 1).Value = $u[0]
    $ws.Cells.Item($row, 2).Value = "image"
    $ws.Cells.Item($row, 3).Value = $u[1]
    $row++
}

# Scroll the view down to the newly added rows and select the last one,
# matching where Excel leaves the cursor after typing in the new data.
$ws.Range("C76").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
